$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 5 -> 4, Wrong mark -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total 125 -> 100, Wrong total -3 -> -6, Max string updated
$ws.Range("B12").Value = 100
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "94 / 112"
